$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.024.32'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '2.054.24'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'248.56"
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = "'0.675"
$ws.Range('E6').Value = '  +3.47%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'54.64"
$ws.Range('E8').Value = '  +15.55%  '
$ws.Range('D9').Value = "'60.66"
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('E10').Value = '  +2.37%  '
$ws.Range('D11').Value = "'0.0785"
$ws.Range('E11').Value = '  +6.07%  '
$ws.Range('E12').Value = '  +6.34%  '
$ws.Range('D13').Value = "'14.87"
$ws.Range('E13').Value = '  +3.13%  '
$ws.Range('D14').Value = '2.354.25'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('E16').Value = '  +3.55%  '
$ws.Range('D17').Value = '2.055.15'
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '36.995.84'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').Value = '0.0₃0933'
$ws.Range('E19').Value = '  +13.30%  '
$ws.Range('D20').Value = "'72.69"
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +7.86%  '
$ws.Range('D22').Value = "'5.33"
$ws.Range('E22').Value = '  +4.20%  '
$ws.Range('D23').Value = "'235.74"
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = "'2.42"
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('D26').Value = "'170.29"
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').Value = "'8.97"
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('D28').Value = "'19.99"
$ws.Range('E28').Value = '  -5.16%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').Value = "'4.57"
$ws.Range('E31').Value = '  +3.44%  '
$ws.Range('D32').Value = "'0.0622"
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('E33').Value = '  +6.06%  '
$ws.Range('D34').Value = "'4.33"
$ws.Range('E34').Value = '  +7.22%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = "'0.0866"
$ws.Range('E36').Value = '  -6.79%  '
$ws.Range('D37').Value = "'2.27"
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('E38').Value = '  -5.70%  '
$ws.Range('D39').Value = "'1.34"
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('D40').Value = "'0.102"
$ws.Range('E40').Value = '  +22.08%  '
$ws.Range('D41').Value = "'17.64"
$ws.Range('E41').Value = '  +11.12%  '
$ws.Range('D42').Value = "'0.0223"
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').Value = "'96.10"
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('D45').Value = "'2.80"
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('D46').Value = "'4.13"
$ws.Range('E46').Value = '  +51.01%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'2.39"
$ws.Range('E47').Value = '  +8.31%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.290.18'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = "'2.92"
$ws.Range('E49').Value = '  +2.98%  '
$ws.Range('B50').Value = 'Gas'
$ws.Range('C50').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D50').Value = "'12.97"
$ws.Range('E50').Value = '  -53.53%  '
$ws.Range('E51').Value = '  +7.03%  '
